# Update workbook to reflect wc_lang ontology renames and related edits.
$wb = $excel.ActiveWorkbook

# --- Submodels sheet: rename "Algorithm" header -> "Framework", value "ssa" -> "stochastic_simulation_algorithm"
$wsSubmodels = $wb.Worksheets.Item("Submodels")
$wsSubmodels.Range("C1").Value = "Framework"
$wsSubmodels.Range("C2").Value = "stochastic_simulation_algorithm"

# --- Compartments sheet: rename geometry/distribution ontology terms
$wsCompartments = $wb.Worksheets.Item("Compartments")
$wsCompartments.Range("C2").Value = "cellular_compartment"
$wsCompartments.Range("D2").Value = "fluid_compartment"
$wsCompartments.Range("E2").Value = "3D_compartment"
$wsCompartments.Range("H2").Value = "normal_distribution"

# --- Initial species concentrations sheet: rename "normal" -> "normal_distribution"
$wsInitConc = $wb.Worksheets.Item("Initial species concentrations")
$wsInitConc.Range("D2").Value = "normal_distribution"
$wsInitConc.Range("D3").Value = "normal_distribution"

# --- Rate laws sheet: clear the "other" Type value
$wsRateLaws = $wb.Worksheets.Item("Rate laws")
$wsRateLaws.Range("E2").Value = ""

# --- Parameters sheet: clear the "other" Type values
$wsParameters = $wb.Worksheets.Item("Parameters")
$wsParameters.Range("C2").Value = ""
$wsParameters.Range("C3").Value = ""
$wsParameters.Range("C4").Value = ""
$wsParameters.Range("C6").Value = ""
$wsParameters.Range("C7").Value = ""

# --- Update sheet selections to match the new active cells
$wsSubmodels.Activate()
$wsSubmodels.Range("C2").Select()

$wsCompartments.Activate()
$wsCompartments.Range("H2").Select()

$wsInitConc.Activate()
$wsInitConc.Range("D2:D3").Select()

$wsRateLaws.Activate()
$wsRateLaws.Range("E2").Select()

# --- Parameters becomes the final active / selected sheet
$wsParameters.Activate()
$wsParameters.Range("C7").Select()
